# Daily attendance processing - 2025-12-24 05:35:25
# Normalize the "Recorded By" (column G) list so that "System" is always
# listed first among the comma-separated recorder names, by swapping the
# first and last entries whenever the cell's recorder list ends with
# ", System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $value = $cell.Text

    if ($value -and $value.EndsWith(", System")) {
        $parts = $value -split ", "
        $last = $parts.Count - 1

        $tmp = $parts[0]
        $parts[0] = $parts[$last]
        $parts[$last] = $tmp

        $cell.Value = $parts -join ", "
    }
}
